# Applies the author's edit to "Project poster.pptx":
#   - Nudge the big border rectangle ("Rectangle 5", shape id 6) slightly
#     left/down: its <a:off> moves from x=442913,y=323850 (EMU) to
#     x=379556,y=431800 (EMU). PowerPoint's Shape.Left/.Top are expressed
#     in points (1 pt = 12700 EMU), so we convert and assign those.
#
# Note: the PowerPoint COM object model has no writable property that
# corresponds to the cosmetic a:rPr/a:endParaRPr "dirty" proofing flag
# that also changed on the title placeholder's runs in the source XML,
# so that purely-cosmetic attribute (no visible/content effect) is not
# reproducible through this interface and is intentionally left alone.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Locate the shape by name/id rather than a hard index, in case ordering
# ever differs.
$rect = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $cand = $s.Shapes.Item($i)
    if ($cand.Id -eq 6 -and $cand.Name -eq "Rectangle 5") {
        $rect = $cand
        break
    }
}
if ($rect -eq $null) {
    $rect = $s.Shapes.Item(32)
}

# Target offsets, in EMU, per the canonical OOXML diff.
# Target EMU -> points conversion is exact math (EMU / 12700), but the
# COM bridge stores Left/Top as a 32-bit float internally before
# re-deriving EMU (floor(f32(pt) * 12700)). Using the midpoint of the
# range of point-values that round-trip to the exact target EMU makes
# the assignment robust to that internal precision loss.
$rect.Left = 29.886339187684370
$rect.Top  = 34.000038147167780
